$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_G"

$values = @(
  71.851851851851862,
  71.851851851851862,
  71.851851851851862,
  71.481481481481481,
  71.851851851851862,
  71.851851851851862,
  73.703703703703709,
  74.074074074074076,
  74.81481481481481,
  74.074074074074076,
  72.592592592592595,
  73.333333333333329,
  74.81481481481481,
  74.81481481481481,
  74.81481481481481,
  74.444444444444443,
  74.444444444444443,
  74.81481481481481,
  73.333333333333329,
  73.703703703703709,
  74.074074074074076,
  70.370370370370367,
  70.370370370370367,
  70,
  73.333333333333329,
  72.962962962962962,
  74.074074074074076,
  75.18518518518519,
  74.81481481481481,
  74.444444444444443,
  70,
  70.370370370370367,
  71.851851851851862,
  70.740740740740733,
  70.740740740740733,
  75.925925925925924,
  69.259259259259252,
  69.629629629629633,
  70.740740740740733,
  74.444444444444443,
  74.444444444444443,
  75.18518518518519,
  74.444444444444443,
  70.740740740740733,
  71.481481481481481,
  70.740740740740733,
  74.81481481481481,
  71.111111111111114
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $values[$i]
}
